$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 45476284
$ws.Range("I40").Value = 18853.166
$ws.Range("J40").Value = 100025200
$ws.Range("K40").Value = 18853.166
$ws.Range("L40").Value = 100025200
$ws.Range("M40").Value = -18678.166
$ws.Range("N40").Value = -100025550
$ws.Range("H64").Value = 24243442
$ws.Range("I64").Value = 8156320
$ws.Range("J64").Value = 41671156
$ws.Range("K64").Value = 8156320
$ws.Range("L64").Value = 41671156
$ws.Range("M64").Value = -8156072
$ws.Range("N64").Value = -41671652
$ws.Range("H67").Value = 24243442
$ws.Range("I67").Value = 8156320
$ws.Range("J67").Value = 41671156
$ws.Range("K67").Value = 8156320
$ws.Range("L67").Value = 41671156
$ws.Range("M67").Value = -8155462
$ws.Range("N67").Value = -41672872
$ws.Range("H132").Value = 124084.664
$ws.Range("I132").Value = 380828.66
$ws.Range("J132").Value = 14051.518
$ws.Range("K132").Value = 1142485.98
$ws.Range("L132").Value = 42154.554
$ws.Range("M132").Value = -1139955.98
$ws.Range("N132").Value = -47214.554
$ws.Range("H137").Value = 3578.186
$ws.Range("I137").Value = 2708.4
$ws.Range("J137").Value = 4044.1428
$ws.Range("K137").Value = 8125.200000000001
$ws.Range("L137").Value = 12132.4284
$ws.Range("M137").Value = -5575.200000000001
$ws.Range("N137").Value = -17232.4284
$ws.Range("H138").Value = 5494.965
$ws.Range("I138").Value = 2519.5833
$ws.Range("J138").Value = 6288.4
$ws.Range("K138").Value = 7558.749899999999
$ws.Range("L138").Value = 18865.2
$ws.Range("M138").Value = -2418.749899999999
$ws.Range("N138").Value = -29145.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 791890.25
$ws.Range("I2").Value = 949444.1
$ws.Range("J2").Value = 4121
$ws.Range("K2").Value = 949444.1
$ws.Range("L2").Value = 4121
$ws.Range("M2").Value = -949331.1
$ws.Range("N2").Value = -4347
$ws.Range("H32").Value = 3638.0144
$ws.Range("I32").Value = 2013.25
$ws.Range("K32").Value = 2013.25
$ws.Range("M32").Value = -1726.25
$ws.Range("H74").Value = 6804.7144
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 6804.7144
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").Value = 6804.7144
$ws.Range("N74").Value = -8552.714400000001
$ws.Range("H77").Value = 6804.7144
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 6804.7144
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").Value = 34023.572
$ws.Range("N77").Value = -42759.572
$ws.Range("H116").Value = 791890.25
$ws.Range("I116").Value = 949444.1
$ws.Range("J116").Value = 4121
$ws.Range("K116").Value = 949444.1
$ws.Range("L116").Value = 4121
$ws.Range("M116").Value = -947150.1
$ws.Range("N116").Value = -8709
$ws.Range("H122").Value = 4590.913
$ws.Range("I122").Value = 3311.0557
$ws.Range("J122").Value = 9198.4
$ws.Range("K122").Value = 9933.167099999999
$ws.Range("L122").Value = 27595.2
$ws.Range("M122").Value = -7483.167099999999
$ws.Range("N122").Value = -32495.2
$ws.Range("H132").Value = 29113.191
$ws.Range("I132").Value = 33376.5
$ws.Range("K132").Value = 100129.5
$ws.Range("M132").Value = -97599.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 791890.25
$ws.Range("I3").Value = 949444.1
$ws.Range("J3").Value = 4121
$ws.Range("K3").Value = 949444.1
$ws.Range("L3").Value = 4121
$ws.Range("M3").Value = -949330.1
$ws.Range("N3").Value = -4349
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2578.52
$ws.Range("I31").Value = 2375.5088
$ws.Range("J31").Value = 2847.628
$ws.Range("K31").Value = 2375.5088
$ws.Range("L31").Value = 2847.628
$ws.Range("M31").Value = -2080.5088
$ws.Range("N31").Value = -3437.628
$ws.Range("H32").Value = 2011
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 2578.52
$ws.Range("I34").Value = 2375.5088
$ws.Range("J34").Value = 2847.628
$ws.Range("K34").Value = 2375.5088
$ws.Range("L34").Value = 2847.628
$ws.Range("M34").Value = -2173.5088
$ws.Range("N34").Value = -3251.628
$ws.Range("H58").Value = 1430288.6
$ws.Range("I58").Value = 1430288.6
$ws.Range("K58").Value = 1430288.6
$ws.Range("M58").Value = -1430085.6
$ws.Range("H132").Value = 25661198
$ws.Range("I132").Value = 33339566
$ws.Range("J132").Value = 66638
$ws.Range("K132").Value = 100018698
$ws.Range("L132").Value = 199914
$ws.Range("M132").Value = -100016168
$ws.Range("N132").Value = -204974
$ws.Range("H134").Value = 3282.75
$ws.Range("I134").Value = 3194.5715
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 9583.7145
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -7048.7145
$ws.Range("N134").Value = -16770
$ws.Range("H136").Value = 1430288.6
$ws.Range("I136").Value = 1430288.6
$ws.Range("K136").Value = 4290865.800000001
$ws.Range("M136").Value = -4288315.800000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 270377.34
$ws.Range("I68").Value = 1692.4667
$ws.Range("J68").Value = 558254
$ws.Range("K68").Value = 5077.4001
$ws.Range("L68").Value = 1674762
$ws.Range("M68").Value = -4266.4001
$ws.Range("N68").Value = -1676384
$ws.Range("H71").Value = 270377.34
$ws.Range("I71").Value = 1692.4667
$ws.Range("J71").Value = 558254
$ws.Range("K71").Value = 15232.2003
$ws.Range("L71").Value = 5024286
$ws.Range("M71").Value = -11176.2003
$ws.Range("N71").Value = -5032398
$ws.Range("H113").Value = 576.8125
$ws.Range("I113").Value = 477.57144
$ws.Range("J113").Value = 654
$ws.Range("K113").Value = 1432.71432
$ws.Range("L113").Value = 1962
$ws.Range("M113").Value = 737.28568
$ws.Range("N113").Value = -6302
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14852.167
$ws.Range("H102").Value = 600543.5
$ws.Range("I102").Value = 1444922.9
$ws.Range("J102").Value = 9477.9
$ws.Range("K102").Value = 1444922.9
$ws.Range("L102").Value = 9477.9
$ws.Range("M102").Value = -1443300.9
$ws.Range("N102").Value = -12721.9
$ws.Range("H132").Value = 3400.394
$ws.Range("I132").Value = 2932.0454
$ws.Range("J132").Value = 4337.091
$ws.Range("K132").Value = 8796.136200000001
$ws.Range("L132").Value = 13011.273
$ws.Range("M132").Value = -6266.136200000001
$ws.Range("N132").Value = -18071.273
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4769.37
$ws.Range("I132").Value = 4251.5557
$ws.Range("J132").Value = 6100.893
$ws.Range("K132").Value = 12754.6671
$ws.Range("L132").Value = 18302.679
$ws.Range("M132").Value = -10224.6671
$ws.Range("N132").Value = -23362.679
$ws.Range("H136").Value = 4887.216
$ws.Range("I136").Value = 4204.6387
$ws.Range("K136").Value = 12613.9161
$ws.Range("M136").Value = -10063.9161
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4318.4287
$ws.Range("I96").Value = 3588.1
$ws.Range("J96").Value = 6144.25
$ws.Range("K96").Value = 3588.1
$ws.Range("L96").Value = 6144.25
$ws.Range("M96").Value = -2215.1
$ws.Range("N96").Value = -8890.25
$ws.Range("H132").Value = 9641251
$ws.Range("I132").Value = 30297.227
$ws.Range("J132").Value = 62501500
$ws.Range("K132").Value = 90891.681
$ws.Range("L132").Value = 187504500
$ws.Range("M132").Value = -88361.681
$ws.Range("N132").Value = -187509560
$ws.Range("H136").Value = 9979.92
$ws.Range("I136").Value = 8000
$ws.Range("J136").Value = 9999.919
$ws.Range("K136").Value = 24000
$ws.Range("L136").Value = 29999.757
$ws.Range("M136").Value = -21450
$ws.Range("N136").Value = -35099.757
